$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 508
$ws.Range("I18").Value = 487.5
$ws.Range("K18").Value = 487.5
$ws.Range("M18").Value = -203.5
$ws.Range("H19").Value = 244.25
$ws.Range("I19").Value = 266.25
$ws.Range("J19").Value = 222.25
$ws.Range("K19").Value = 266.25
$ws.Range("L19").Value = 222.25
$ws.Range("M19").Value = -91.25
$ws.Range("N19").Value = -572.25
$ws.Range("H33").Value = 169.03572
$ws.Range("I33").Value = 172.03847
$ws.Range("K33").Value = 172.03847
$ws.Range("M33").Value = 56.96153000000001
$ws.Range("H86").Value = 9207.385
$ws.Range("I86").Value = 1616.6666
$ws.Range("K86").Value = 1616.6666
$ws.Range("M86").Value = -493.6666
$ws.Range("H89").Value = 9207.385
$ws.Range("I89").Value = 1616.6666
$ws.Range("K89").Value = 8083.333000000001
$ws.Range("M89").Value = -2467.333000000001
$ws.Range("H100").Value = 2244.111
$ws.Range("I100").Value = 1599.4
$ws.Range("J100").Value = 3050
$ws.Range("K100").Value = 1599.4
$ws.Range("L100").Value = 3050
$ws.Range("M100").Value = -1058.4
$ws.Range("N100").Value = -4132
$ws.Range("H103").Value = 83333580
$ws.Range("J103").Value = 400
$ws.Range("L103").Value = 1200
$ws.Range("N103").Value = -2372
$ws.Range("H129").Value = 1278.475
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 1301
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 3903
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -13903
$ws.Range("H135").Value = 12198078
$ws.Range("I135").Value = 570.25714
$ws.Range("J135").Value = 83350210
$ws.Range("K135").Value = 5132.31426
$ws.Range("L135").Value = 750151890
$ws.Range("M135").Value = -2597.31426
$ws.Range("N135").Value = -750156960
$ws.Range("H138").Value = 135907.42
$ws.Range("J138").Value = 151532.56
$ws.Range("L138").Value = 454597.68
$ws.Range("N138").Value = -464877.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1149.375
$ws.Range("I2").Value = 1149.1666
$ws.Range("K2").Value = 1149.1666
$ws.Range("M2").Value = -1036.1666
$ws.Range("H61").Value = 2915.8206
$ws.Range("I61").Value = 2717.1667
$ws.Range("K61").Value = 2717.1667
$ws.Range("M61").Value = -2505.1667
$ws.Range("H74").Value = 22223262
$ws.Range("I74").Value = 25000554
$ws.Range("J74").Value = 4920
$ws.Range("K74").Value = 25000554
$ws.Range("L74").Value = 4920
$ws.Range("M74").Value = -24999680
$ws.Range("N74").Value = -6668
$ws.Range("H77").Value = 22223262
$ws.Range("I77").Value = 25000554
$ws.Range("J77").Value = 4920
$ws.Range("K77").Value = 125002770
$ws.Range("L77").Value = 24600
$ws.Range("M77").Value = -124998402
$ws.Range("N77").Value = -33336
$ws.Range("H97").Value = 533.25
$ws.Range("I97").Value = 536.2727
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 536.2727
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -40.27269999999999
$ws.Range("N97").Value = -1492
$ws.Range("H102").Value = 1388.75
$ws.Range("I102").Value = 1280.8823
$ws.Range("K102").Value = 1280.8823
$ws.Range("M102").Value = 341.1177
$ws.Range("H116").Value = 1149.375
$ws.Range("I116").Value = 1149.1666
$ws.Range("K116").Value = 1149.1666
$ws.Range("M116").Value = 1144.8334
$ws.Range("H122").Value = 2066.8965
$ws.Range("I122").Value = 1858.4
$ws.Range("K122").Value = 5575.200000000001
$ws.Range("M122").Value = -3125.200000000001
$ws.Range("H132").Value = 12252.714
$ws.Range("I132").Value = 1705.9269
$ws.Range("K132").Value = 5117.780699999999
$ws.Range("M132").Value = -2587.780699999999
$ws.Range("H136").Value = 2915.8206
$ws.Range("I136").Value = 2717.1667
$ws.Range("K136").Value = 8151.500100000001
$ws.Range("M136").Value = -5601.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1149.375
$ws.Range("I3").Value = 1149.1666
$ws.Range("K3").Value = 1149.1666
$ws.Range("M3").Value = -1035.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 172.7619
$ws.Range("I22").Value = 156.76923
$ws.Range("J22").Value = 198.75
$ws.Range("K22").Value = 156.76923
$ws.Range("L22").Value = 198.75
$ws.Range("M22").Value = 193.23077
$ws.Range("N22").Value = -898.75
$ws.Range("H31").Value = 3719
$ws.Range("I31").Value = 1796.5714
$ws.Range("J31").Value = 5872.12
$ws.Range("K31").Value = 1796.5714
$ws.Range("L31").Value = 5872.12
$ws.Range("M31").Value = -1501.5714
$ws.Range("N31").Value = -6462.12
$ws.Range("H34").Value = 3719
$ws.Range("I34").Value = 1796.5714
$ws.Range("J34").Value = 5872.12
$ws.Range("K34").Value = 1796.5714
$ws.Range("L34").Value = 5872.12
$ws.Range("M34").Value = -1594.5714
$ws.Range("N34").Value = -6276.12
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H60").Value = 20125
$ws.Range("J60").Value = 20142.857
$ws.Range("L60").Value = 20142.857
$ws.Range("N60").Value = -21164.857
$ws.Range("H68").Value = 53437.332
$ws.Range("J68").Value = 53437.332
$ws.Range("L68").Value = 53437.332
$ws.Range("N68").Value = -54935.332
$ws.Range("H71").Value = 53437.332
$ws.Range("J71").Value = 53437.332
$ws.Range("L71").Value = 160311.996
$ws.Range("N71").Value = -167799.996
$ws.Range("H74").Value = 41485.582
$ws.Range("J74").Value = 41485.582
$ws.Range("L74").Value = 41485.582
$ws.Range("N74").Value = -43233.582
$ws.Range("H77").Value = 41485.582
$ws.Range("J77").Value = 41485.582
$ws.Range("L77").Value = 124456.746
$ws.Range("N77").Value = -133192.746
$ws.Range("H95").Value = 26000
$ws.Range("J95").Value = 26000
$ws.Range("L95").Value = 26000
$ws.Range("N95").Value = -31492
$ws.Range("H109").Value = 123772370
$ws.Range("J109").Value = 123772370
$ws.Range("L109").Value = 123772370
$ws.Range("N109").Value = -123774450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 618.91
$ws.Range("J131").Value = 748.5753
$ws.Range("L131").Value = 2245.7259
$ws.Range("N131").Value = -12325.7259
$ws.Range("H137").Value = 11498242
$ws.Range("I137").Value = 1012.8571
$ws.Range("J137").Value = 15156452
$ws.Range("K137").Value = 3038.5713
$ws.Range("L137").Value = 45469356
$ws.Range("M137").Value = 2061.4287
$ws.Range("N137").Value = -45479556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 158.52942
$ws.Range("I55").Value = 134.33333
$ws.Range("J55").Value = 171.72728
$ws.Range("K55").Value = 134.33333
$ws.Range("L55").Value = 171.72728
$ws.Range("M55").Value = 38.66667000000001
$ws.Range("N55").Value = -517.7272800000001
$ws.Range("H68").Value = 2670.7144
$ws.Range("I68").Value = 1700
$ws.Range("K68").Value = 1700
$ws.Range("M68").Value = -951
$ws.Range("H71").Value = 2670.7144
$ws.Range("I71").Value = 1700
$ws.Range("K71").Value = 8500
$ws.Range("M71").Value = -4756

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 277.42856
$ws.Range("I81").Value = 290
$ws.Range("J81").Value = 202
$ws.Range("K81").Value = 580
$ws.Range("L81").Value = 404
$ws.Range("M81").Value = 481
$ws.Range("N81").Value = -2526
$ws.Range("H84").Value = 277.42856
$ws.Range("I84").Value = 290
$ws.Range("J84").Value = 202
$ws.Range("K84").Value = 2900
$ws.Range("L84").Value = 2020
$ws.Range("M84").Value = 2404
$ws.Range("N84").Value = -12628
